$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the set_voltage column (G) from 55.6 VDC to 55.2 VDC (absorption voltage)
# for the rows that were using the old value.
$rows = @(3, 4, 5, 6, 7, 8, 24, 25, 26)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = 55.2
}

# Move the active cell selection to H16
$ws.Range("H16").Select()
